# Add files via upload
# Inserts 8 new "AL-AQEEL GOLD" price rows into Sheet1, right before the
# existing "AL-AQEEL SILVER" block (old row 27), matching the new product
# lines: GOLD SATUAN / PAKET ISI 3 / ISI 5 / ISI 7, each offered in two
# sizes (A5 and A7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string labels used by the inserted rows.
$goldSatuan = "AL-AQEEL GOLD SATUAN"
$goldIsi3   = "AL-AQEEL GOLD PAKET ISI 3"
$goldIsi5   = "AL-AQEEL GOLD PAKET ISI 5"
$goldIsi7   = "AL-AQEEL GOLD PAKET ISI 7"

$hvs = "HVS"
$a5  = "A5 (14,8X21CM)"
$a7  = "A7 (7,4X10,5CM)"

# Row data (B, C, D, E) for the 8 new rows, in the order they must appear,
# all inserted directly above the old row 27 ("AL-AQEEL SILVER" B5 size).
$newRows = @(
    @($goldSatuan, $hvs, $a5, 21000),
    @($goldIsi3,   $hvs, $a5, 63000),
    @($goldIsi5,   $hvs, $a5, 105000),
    @($goldIsi7,   $hvs, $a5, 147000),
    @($goldSatuan, $hvs, $a7, 8500),
    @($goldIsi3,   $hvs, $a7, 25500),
    @($goldIsi5,   $hvs, $a7, 42500),
    @($goldIsi7,   $hvs, $a7, 59500)
)

$insertAt = 27

foreach ($row in $newRows) {
    # Push everything at/after row 27 down by one, each time re-creating
    # the row directly above the "AL-AQEEL SILVER" block.
    $ws.Rows.Item($insertAt).EntireRow.Insert()

    $rowRange = $ws.Range("A" + $insertAt + ":E" + $insertAt)
    $rowRange.Borders.Item(1).LineStyle = 1
    $rowRange.Borders.Item(2).LineStyle = 1
    $rowRange.Borders.Item(3).LineStyle = 1
    $rowRange.Borders.Item(4).LineStyle = 1

    # Column A (the running "No." counter) is intentionally left blank for
    # these rows, matching the other multi-size/variant groups in the sheet.
    $ws.Cells.Item($insertAt, 2).Value = $row[0]
    $ws.Cells.Item($insertAt, 3).Value = $row[1]
    $ws.Cells.Item($insertAt, 4).Value = $row[2]
    $ws.Cells.Item($insertAt, 5).Value = $row[3]

    $insertAt = $insertAt + 1
}

# Restore the view to roughly where the author left it after the edit.
$ws.Range("A25").Select()
